$d = $word.ActiveDocument

$d.Content.Find.Execute("Registro de tutor almacenado.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Registro de matrícula almacenada.", 2)
